$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.9951798856757821
$ws.Range("E2").Value = 0.9951798856757821

$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1

$ws.Range("D4").Value = 0.00061381034650537
$ws.Range("E4").Value = 0.00061381034650537

$ws.Range("D5").Value = 0.000009737410624520885
$ws.Range("E5").Value = 0.000009737410624520885

$ws.Range("D6").Value = 0.1727573487936249
$ws.Range("E6").Value = 0.1727573487936249

$ws.Range("D7").Value = 0.9999999999999873
$ws.Range("E7").Value = 0.00000000000001265654248072678

$ws.Range("D8").Value = 0.9999999997397482
$ws.Range("E8").Value = 0.0000000002602518200944814

$ws.Range("D9").Value = 0.9993778928338607
$ws.Range("E9").Value = 0.0006221071661393074

$ws.Range("D11").Value = 0.9999999999999851
$ws.Range("E11").Value = 0.0000000000000148769885299771
$ws.Range("F11").Value = 4.321595668792725
